# Generate Report for Handback
# Update timestamp values for the bc329aef-741a-4619-b73e-cbabdb3096fc row
# across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for bc329aef row (row 4)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-31 00:48:38"

# zh-cn sheet: "Correspond Handoff Datetime" (H) and "Correspond Handback DateTime" (K) for bc329aef row (row 4)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-31 00:48:33"
$wsZhCn.Range("K4").Value = "2016-08-31 00:49:09"

# de-de sheet: "Correspond Handback DateTime" (K) for bc329aef row (row 4)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K4").Value = "2016-08-31 00:49:17"
